# Anchor regex patterns with ^ and $ in the "(REGEX: ...)" annotations.
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "(REGEX: fr(\.[\w-]+){3,4})", $false, $false, $false, $false, $false,
    $true, 1, $false, "(REGEX: ^fr(\.[\w-]+){3,4}$)", 2)

$d.Content.Find.Execute(
    "(REGEX: ([\w-]+\.){3,4}resource(\.[\w-]+){1,2})", $false, $false, $false, $false, $false,
    $true, 1, $false, "(REGEX: ^([\w-]+\.){3,4}resource(\.[\w-]+){1,2}$)", 2)

$d.Content.Find.Execute(
    "(REGEX: ([\w-]+\.){3,4}request(\.[\w-]+){1,2})", $false, $false, $false, $false, $false,
    $true, 1, $false, "(REGEX: ^([\w-]+\.){3,4}request(\.[\w-]+){1,2}$)", 2)

$d.Content.Find.Execute(
    "(REGEX: [0-9]{5})", $false, $false, $false, $false, $false,
    $true, 1, $false, "(REGEX:^[0-9]{5}$)", 2)
